$d = $word.ActiveDocument

# Update the date header (unique text, safe via Find/Replace).
$d.Content.Find.Execute("2025-08-28 Thursday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-08-29 Friday", 2)

# Update the division-problem table. Addressed by (row, column) rather than
# Find/Replace because several old/new values collide across cells
# (e.g. "77÷5=" is simultaneously an old value in one cell and the new
# value of another, and "43÷6=" is an old value in two different cells),
# so direct cell assignment avoids any re-matching issues.
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "75÷6="
$t.Cell(1, 2).Range.Text = "66÷5="
$t.Cell(1, 3).Range.Text = "12÷2="
$t.Cell(1, 4).Range.Text = "94÷6="
$t.Cell(1, 5).Range.Text = "82÷3="

$t.Cell(5, 1).Range.Text = "73÷7="
$t.Cell(5, 2).Range.Text = "52÷8="
$t.Cell(5, 3).Range.Text = "87÷3="
$t.Cell(5, 4).Range.Text = "55÷9="
$t.Cell(5, 5).Range.Text = "77÷5="

$t.Cell(9, 1).Range.Text = "82÷5="
$t.Cell(9, 2).Range.Text = "49÷7="
$t.Cell(9, 3).Range.Text = "61÷5="
$t.Cell(9, 4).Range.Text = "97÷9="
$t.Cell(9, 5).Range.Text = "50÷4="

$t.Cell(13, 1).Range.Text = "88÷3="
$t.Cell(13, 2).Range.Text = "18÷4="
$t.Cell(13, 3).Range.Text = "34÷8="
$t.Cell(13, 4).Range.Text = "39÷5="
$t.Cell(13, 5).Range.Text = "12÷9="

$t.Cell(17, 1).Range.Text = "76÷8="
$t.Cell(17, 2).Range.Text = "72÷6="
$t.Cell(17, 3).Range.Text = "17÷7="
$t.Cell(17, 4).Range.Text = "19÷6="
$t.Cell(17, 5).Range.Text = "83÷6="

Write-Host "Done."
